$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLPI Input")

# Capture the current (pre-edit) pixel widths of the columns that are about to
# shift, so the same visual widths can be re-applied to their new homes once
# the two new columns have been inserted.
$wFabricanteTelefone = $ws.Columns("N").Width   # -> becomes "Marca do Celular" (O)
$wModeloCelular      = $ws.Columns("O").Width   # -> becomes "Modelo do Celular" (P)
$wImeiCelular        = $ws.Columns("P").Width   # -> becomes "IMEI do Celular" (Q)

# Insert two new columns right before the "Fabricante telefone" column (N) to
# hold the two new phone attributes.
$ws.Columns("N:O").Insert()

# New headers for the inserted columns.
$ws.Range("N1").Value = "Tipo de Celular"
$ws.Range("O1").Value = "Marca do Celular"

# Match the header cell shading used by the sibling phone-info headers
# (style was inherited from the column to the left on insert, but "Marca do
# Celular" should look like the other plain headers, e.g. "Entidade A").
$ws.Range("D1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

# The old "Fabricante telefone" and "Tipo de Telefone" columns are now
# redundant (replaced conceptually by "Tipo de Celular" / "Marca do
# Celular"), so remove them. After the insert they sit at P ("Fabricante
# telefone") and S ("Tipo de Telefone"); delete from right to left so the
# earlier deletion doesn't shift the later column reference.
$ws.Columns("S:S").Delete()
$ws.Columns("P:P").Delete()

# Restore the visual widths for the columns that shifted right by one slot,
# and give the brand-new "Tipo de Celular" column the same plain width used
# by its neighbouring custom-width (non-autofit) columns.
$ws.Columns("N").Width = $wImeiCelular
$ws.Columns("O").Width = $wFabricanteTelefone
$ws.Columns("P").Width = $wModeloCelular
$ws.Columns("Q").Width = $wImeiCelular

# Move the selection, matching where editing left off.
$ws.Range("H36").Select() | Out-Null
